$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update activation date (shared by B8/C8 and B15/C15 since they reused the same text)
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "01/01/2023"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2023"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2023"

# 2. Fill in the English "Objectives" text (row 11) - new cells
$objectivesText = "To present experimental techniques for the characterization of electrical, magnetic and thermal  properties of materials."
$ws.Range("B11").Value = $objectivesText
$ws.Range("C11").Value = $objectivesText

# 3. Fill in the "Short syllabus" text (row 14) - reuses same text as row 11
$ws.Range("B14").Value = $objectivesText
$ws.Range("C14").Value = $objectivesText

# 4. Fill in the "Syllabus" text (row 16) - new cells
$syllabusText = "histerese de materiais magnéticos macios. Medidas de magnetostricção.Propriedades térmicas dos materiais:  expansão térmica.Electrical properties: electrical conductivity in pure metals, metallic alloys and semiconductors, and superconductors; Hall Effect; Ohm's Law and dependence on temperature.Magnetic properties: magnetic susceptibility and c.c. magnetization. Hysteresis curves of soft magnetic materials. Magnetostriction measurements.Thermal properties of materials: thermal expansion."
$ws.Range("B16").Value = $syllabusText
$ws.Range("C16").Value = $syllabusText

# 5. Update the "Norma de recuperação" text (row 20)
$ws.Range("B20").Value = "Média aritmética das notas dos relatórios de cada experimento"
$ws.Range("C20").Value = "Média aritmética das notas dos relatórios de cada experimento"
